# Update cryptocurrency price/volume data per latest symbol-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'309.88"
$ws.Range('E2').Value = "'-3.60%"
$ws.Range('D3').Value = "'49.43"
$ws.Range('E3').Value = "'1.00%"
$ws.Range('D4').Value = "'5.172"
$ws.Range('E4').Value = "'-1.58%"
$ws.Range('D5').Value = "'0.07778"
$ws.Range('E5').Value = "'-3.84%"
$ws.Range('D6').Value = "'4.519"
$ws.Range('E6').Value = "'-2.10%"
$ws.Range('D7').Value = "'1.381"
$ws.Range('E7').Value = "'14.58%"
$ws.Range('D8').Value = "'1.554"
$ws.Range('E8').Value = "'-6.36%"
$ws.Range('D9').Value = "'0.1231"
$ws.Range('E9').Value = "'-6.09%"
$ws.Range('D10').Value = "'0.2008"
$ws.Range('E10').Value = "'3.14%"
$ws.Range('D11').Value = "'0.09510"
$ws.Range('E11').Value = "'-0.45%"
$ws.Range('D12').Value = "'0.04685"
$ws.Range('E12').Value = "'4.79%"
$ws.Range('E13').Value = "'-0.15%"
$ws.Range('D14').Value = "'0.001270"
$ws.Range('E14').Value = "'-4.25%"
$ws.Range('D15').Value = "'0.04174"
$ws.Range('E15').Value = "'-3.18%"
$ws.Range('D16').Value = "'0.005802"
$ws.Range('E16').Value = "'-2.36%"
$ws.Range('E17').Value = "'2,016.87%"
$ws.Range('D18').Value = "'3.337"
$ws.Range('E18').Value = "'-0.64%"
$ws.Range('D19').Value = "'2.239"
$ws.Range('E19').Value = "'-8.31%"
$ws.Range('D20').Value = "'0.3449"
$ws.Range('E20').Value = "'1.71%"
$ws.Range('D21').Value = "'7.935"
$ws.Range('E21').Value = "'-3.56%"
$ws.Range('D22').Value = "'0.1341"
$ws.Range('E22').Value = "'-5.11%"
$ws.Range('E23').Value = "'3.91%"
$ws.Range('D24').Value = "'0.001271"
$ws.Range('E24').Value = "'-2.87%"
$ws.Range('D25').Value = "'0.004045"
$ws.Range('E25').Value = "'-4.72%"
$ws.Range('E26').Value = "'-0.21%"
$ws.Range('D38').Value = "'0.02609"
$ws.Range('E38').Value = "'-2.87%"
$ws.Range('D39').Value = "'0.05868"
$ws.Range('E39').Value = "'5.20%"
$ws.Range('D40').Value = "'0.01073"
$ws.Range('E40').Value = "'69.88%"
$ws.Range('D41').Value = "'0.007941"
$ws.Range('E41').Value = "'3.22%"
$ws.Range('D42').Value = "'0.1425"
$ws.Range('E42').Value = "'-0.91%"
$ws.Range('D43').Value = "'0.008426"
$ws.Range('E43').Value = "'9.28%"
$ws.Range('D44').Value = "'0.008314"
$ws.Range('E44').Value = "'2.63%"
$ws.Range('D45').Value = "'0.3404"
$ws.Range('E45').Value = "'6.59%"
$ws.Range('D46').Value = "'0.00007024"
$ws.Range('E46').Value = "'-0.05%"
$ws.Range('E47').Value = "'-0.20%"
$ws.Range('D48').Value = "'0.05235"
$ws.Range('E48').Value = "'-14.48%"
$ws.Range('D49').Value = "'0.002621"
$ws.Range('E49').Value = "'-34.63%"
$ws.Range('D50').Value = "'0.00002101"
$ws.Range('E50').Value = "'-0.20%"
$ws.Range('D51').Value = "'0.0002001"
$ws.Range('E51').Value = "'-0.20%"
